{"js": "// Apply text replacements to update the worksheet date and division problems.\nconst replacements = [\n  [\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"],\n  [\"73\u00f78=9, 1\", \"86\u00f73=28, 2\"],\n  [\"81\u00f73=27, 0\", \"55\u00f75=11, 0\"],\n  [\"66\u00f75=13, 1\", \"14\u00f73=4, 2\"],\n  [\"99\u00f75=19, 4\", \"86\u00f73=28, 2\"],\n  [\"77\u00f74=19, 1\", \"28\u00f79=3, 1\"],\n  [\"63\u00f73=21, 0\", \"78\u00f73=26, 0\"],\n  [\"45\u00f72=22, 1\", \"74\u00f72=37, 0\"],\n  [\"99\u00f74=24, 3\", \"25\u00f73=8, 1\"],\n  [\"40\u00f76=6, 4\", \"78\u00f77=11, 1\"],\n  [\"32\u00f79=3, 5\", \"37\u00f79=4, 1\"],\n  [\"46\u00f72=23, 0\", \"25\u00f79=2, 7\"],\n  [\"38\u00f75=7, 3\", \"99\u00f79=11, 0\"],\n  [\"79\u00f74=19, 3\", \"64\u00f76=10, 4\"],\n  [\"90\u00f72=45, 0\", \"24\u00f77=3, 3\"],\n  [\"88\u00f77=12, 4\", \"91\u00f78=11, 3\"],\n  [\"13\u00f75=2, 3\", \"87\u00f74=21, 3\"],\n  [\"95\u00f79=10, 5\", \"75\u00f74=18, 3\"],\n  [\"48\u00f73=16, 0\", \"66\u00f78=8, 2\"],\n  [\"96\u00f73=32, 0\", \"43\u00f72=21, 1\"],\n  [\"44\u00f79=4, 8\", \"50\u00f72=25, 0\"],\n  [\"82\u00f78=10, 2\", \"42\u00f77=6, 0\"],\n  [\"44\u00f72=22, 0\", \"16\u00f77=2, 2\"],\n  [\"51\u00f78=6, 3\", \"21\u00f78=2, 5\"],\n  [\"50\u00f79=5, 5\", \"95\u00f73=31, 2\"],\n  [\"12\u00f73=4, 0\", \"50\u00f75=10, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "# Update the worksheet date and each division-problem cell to match the new day's values.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n\nReplace-Text \"2025-09-22 Monday\" \"2025-09-23 Tuesday\"\nReplace-Text \"73\u00f78=9, 1\" \"86\u00f73=28, 2\"\nReplace-Text \"81\u00f73=27, 0\" \"55\u00f75=11, 0\"\nReplace-Text \"66\u00f75=13, 1\" \"14\u00f73=4, 2\"\nReplace-Text \"99\u00f75=19, 4\" \"86\u00f73=28, 2\"\nReplace-Text \"77\u00f74=19, 1\" \"28\u00f79=3, 1\"\nReplace-Text \"63\u00f73=21, 0\" \"78\u00f73=26, 0\"\nReplace-Text \"45\u00f72=22, 1\" \"74\u00f72=37, 0\"\nReplace-Text \"99\u00f74=24, 3\" \"25\u00f73=8, 1\"\nReplace-Text \"40\u00f76=6, 4\" \"78\u00f77=11, 1\"\nReplace-Text \"32\u00f79=3, 5\" \"37\u00f79=4, 1\"\nReplace-Text \"46\u00f72=23, 0\" \"25\u00f79=2, 7\"\nReplace-Text \"38\u00f75=7, 3\" \"99\u00f79=11, 0\"\nReplace-Text \"79\u00f74=19, 3\" \"64\u00f76=10, 4\"\nReplace-Text \"90\u00f72=45, 0\" \"24\u00f77=3, 3\"\nReplace-Text \"88\u00f77=12, 4\" \"91\u00f78=11, 3\"\nReplace-Text \"13\u00f75=2, 3\" \"87\u00f74=21, 3\"\nReplace-Text \"95\u00f79=10, 5\" \"75\u00f74=18, 3\"\nReplace-Text \"48\u00f73=16, 0\" \"66\u00f78=8, 2\"\nReplace-Text \"96\u00f73=32, 0\" \"43\u00f72=21, 1\"\nReplace-Text \"44\u00f79=4, 8\" \"50\u00f72=25, 0\"\nReplace-Text \"82\u00f78=10, 2\" \"42\u00f77=6, 0\"\nReplace-Text \"44\u00f72=22, 0\" \"16\u00f77=2, 2\"\nReplace-Text \"51\u00f78=6, 3\" \"21\u00f78=2, 5\"\nReplace-Text \"50\u00f79=5, 5\" \"95\u00f73=31, 2\"\nReplace-Text \"12\u00f73=4, 0\" \"50\u00f75=10, 0\"\n"}
